# Add a third dataset ("cr" = customer reviews) as a new worksheet after "trec".
$wb = $excel.ActiveWorkbook

# Match the "before" state's selection on "trec" (it becomes the non-active sheet
# once "cr" is added, and its selection moves to cover the data table A2:C10).
$trec = $wb.Worksheets.Item("trec")
$trec.Activate()
$trec.Range("A2:C10").Select()

# Insert the new sheet immediately after the last existing sheet ("trec").
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$newSheet.Name = "cr"

# Header row (reuses the same shared strings as the other two sheets).
$headers = @("% dataset", "Regular", "Augmentation")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(2, $i + 1).Value = $headers[$i]
}

# Data rows: % dataset, Regular, Augmentation accuracy for the customer-reviews set.
$data = @(
    @(0.001, 0.53, 0.41),
    @(0.003, 0.41, 0.47),
    @(0.01, 0.49, 0.68),
    @(0.05, 0.59, 0.71),
    @(0.1, 0.65, 0.76),
    @(0.25, 0.69, 0.77),
    @(0.5, 0.74, 0.78),
    @(1, 0.74, 0.66)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt 3; $c++) {
        $newSheet.Cells.Item($r + 3, $c + 1).Value = $data[$r][$c]
    }
}

# Final selection/active state on the new sheet, matching the target workbook.
$newSheet.Range("B14").Select()
